$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0

# Row 3 (status)
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("J3").Value = 1

# Row 4 (name)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0

# Row 5 (address)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0

# Row 6 (latitude)
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0

# Row 7 (longitude)
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0

# Row 8 (region1)
$ws.Range("J8").Value = 1

# Row 9 (type_id)
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0

# Row 10 (type_name)
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0

# Row 12 (product1)
$ws.Range("K12").Value = 0

# Row 13 (product2)
$ws.Range("K13").Value = 0

# Row 14 (product3)
$ws.Range("K14").Value = 0

# Row 15 (service1)
$ws.Range("K15").Value = 0

# Row 16 (service2)
$ws.Range("K16").Value = 0

# Row 17 (service3)
$ws.Range("K17").Value = 0

# Row 18 (service4)
$ws.Range("J18").Value = 1

# Row 19 (service5)
$ws.Range("J19").Value = 1

# Row 20 (region2)
$ws.Range("J20").Value = 1

# Row 21 (region3)
$ws.Range("J21").Value = 1
